# Bug-fix commit: highlight (yellow) the text of three requirement
# paragraphs so they stand out as still-needing-work items.
#   1. "Si chiede dunque di creare delle API che consentano di:"
#   2. "Scalare il credito di un utente ... tempi stabiliti"
#   3. "Prevedere una rotta per l'utente con ruolo admin ... mediante JWT)."
#
# wdYellow = 7 for Range.HighlightColorIndex (-> <w:highlight w:val="yellow"/>)

$d = $word.ActiveDocument
$wdYellow = 7

$matchCount = 0

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -like "Si chiede dunque di creare delle API*") {
        # Single-run paragraph: plain run-level highlight is all that is
        # needed (the paragraph mark itself keeps no highlight).
        $p.Range.HighlightColorIndex = $wdYellow
        $matchCount = $matchCount + 1
    }
    elseif ($t -like "Scalare il credito di u*tempi stabiliti*") {
        # This paragraph's selection, when originally edited, also covered
        # the paragraph mark (e.g. a triple-click selection), so Word
        # additionally recorded the highlight on the paragraph mark's own
        # run properties (w:pPr/w:rPr) - something the simple
        # Range.HighlightColorIndex setter does not reproduce for the
        # paragraph mark. Rebuild the paragraph's OOXML explicitly (adding
        # <w:highlight> to w:pPr/w:rPr and to every run's w:rPr) while
        # keeping every existing paragraph/run identity attribute
        # (w14:paraId, w14:textId, w:rsidR, ...) untouched.
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="62A80B6C" w14:textId="26B707D9" w:rsidR="004152B0" w:rsidRPr="000B6917" w:rsidRDefault="004152B0" w:rsidP="007B0472"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r w:rsidRPr="000B6917"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Scalare il credito di u</w:t></w:r><w:r w:rsidR="00B14C81" w:rsidRPr="000B6917"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>n utente al momento della &#8220;prenotazione&#8221; dello slot. Restituire il credito se viene effettuata la cancellazione entro i tempi stabiliti</w:t></w:r></w:p>'
        $p.Range.InsertXML($xml)
        $matchCount = $matchCount + 1
    }
    elseif ($t -like "Prevedere una rotta per l*utente con ruolo admin*") {
        # Three runs, highlight each one individually (no paragraph-mark
        # highlight here - that paragraph's selection stopped at the
        # closing parenthesis/JWT text).
        $p.Range.HighlightColorIndex = $wdYellow
        $matchCount = $matchCount + 1
    }
}

if ($matchCount -ne 3) {
    throw "expected to highlight 3 paragraphs, matched $matchCount"
}
